# Scheduled runner update: refresh market-price / profit figures across the
# Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) for the rows
# whose underlying item prices changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3269016.5
$ws.Range("J112").Value = 3269016.5
$ws.Range("L112").Value = 9807049.5
$ws.Range("N112").Value = -9809265.5

$ws.Range("H129").Value = 501343.6
$ws.Range("J129").Value = 626605.1
$ws.Range("L129").Value = 1879815.3
$ws.Range("N129").Value = -1889815.3

$ws.Range("H132").Value = 2204.282
$ws.Range("I132").Value = 2312.8918
$ws.Range("K132").Value = 6938.6754
$ws.Range("M132").Value = -4408.6754

$ws.Range("H137").Value = 2057.6365
$ws.Range("I137").Value = 1877.5
$ws.Range("K137").Value = 5632.5
$ws.Range("M137").Value = -3082.5

$ws.Range("H141").Value = 1024.5686
$ws.Range("J141").Value = 3699.75
$ws.Range("L141").Value = 11099.25
$ws.Range("N141").Value = -21459.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H6").Value = 50000000
$ws.Range("I6").Value = 50000000
$ws.Range("K6").Value = 50000000
$ws.Range("M6").Value = -49999827

$ws.Range("H9").Value = 10009
$ws.Range("J9").Value = 10009
$ws.Range("L9").Value = 10009
$ws.Range("N9").Value = -10349

$ws.Range("H20").Value = 10009
$ws.Range("J20").Value = 10009
$ws.Range("L20").Value = 10009
$ws.Range("N20").Value = -10549

$ws.Range("H32").Value = 4528.045
$ws.Range("I32").Value = 4384.213
$ws.Range("K32").Value = 4384.213
$ws.Range("M32").Value = -4097.213

$ws.Range("H37").Value = 29895
$ws.Range("J37").Value = 29895
$ws.Range("L37").Value = 29895
$ws.Range("N37").Value = -30441

$ws.Range("H44").Value = 26000
$ws.Range("J44").Value = 26000
$ws.Range("L44").Value = 26000
$ws.Range("N44").Value = -26976

$ws.Range("H54").Value = 17000
$ws.Range("J54").Value = 17000
$ws.Range("L54").Value = 17000
$ws.Range("N54").Value = -18538

$ws.Range("H74").Value = 111112110
$ws.Range("I74").Value = 125000980
$ws.Range("J74").Value = 1200
$ws.Range("K74").Value = 125000980
$ws.Range("L74").Value = 1200
$ws.Range("M74").Value = -125000106
$ws.Range("N74").Value = -2948

$ws.Range("H77").Value = 111112110
$ws.Range("I77").Value = 125000980
$ws.Range("J77").Value = 1200
$ws.Range("K77").Value = 625004900
$ws.Range("L77").Value = 6000
$ws.Range("M77").Value = -625000532
$ws.Range("N77").Value = -14736

$ws.Range("H97").Value = 1259.3572
$ws.Range("I97").Value = 1308.5454
$ws.Range("J97").Value = 1079
$ws.Range("K97").Value = 1308.5454
$ws.Range("L97").Value = 1079
$ws.Range("M97").Value = -812.5454
$ws.Range("N97").Value = -2071

$ws.Range("H102").Value = 1148.1666
$ws.Range("I102").Value = 847.25
$ws.Range("J102").Value = 1750
$ws.Range("K102").Value = 847.25
$ws.Range("L102").Value = 1750
$ws.Range("M102").Value = 774.75
$ws.Range("N102").Value = -4994

$ws.Range("H132").Value = 27182.61
$ws.Range("I132").Value = 1526.7609
$ws.Range("K132").Value = 4580.2827
$ws.Range("M132").Value = -2050.2827

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1321.65
$ws.Range("I86").Value = 1180.1875
$ws.Range("K86").Value = 1180.1875
$ws.Range("M86").Value = -57.1875

$ws.Range("H89").Value = 1321.65
$ws.Range("I89").Value = 1180.1875
$ws.Range("K89").Value = 5900.9375
$ws.Range("M89").Value = -284.9375

$ws.Range("H99").Value = 2042
$ws.Range("I99").Value = 2052.5
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 2052.5
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -554.5
$ws.Range("N99").Value = -4996

$ws.Range("H129").Value = 49797.4
$ws.Range("J129").Value = 49797.4
$ws.Range("L129").Value = 49797.4
$ws.Range("N129").Value = -59797.4

$ws.Range("H134").Value = 5210.2173
$ws.Range("I134").Value = 5716.8945
$ws.Range("K134").Value = 17150.6835
$ws.Range("M134").Value = -14615.6835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 839.2727
$ws.Range("I16").Value = 484.4
$ws.Range("J16").Value = 1135
$ws.Range("K16").Value = 484.4
$ws.Range("L16").Value = 1135
$ws.Range("M16").Value = -197.4
$ws.Range("N16").Value = -1709

$ws.Range("H31").Value = 11717.195
$ws.Range("I31").Value = 14324.667
$ws.Range("K31").Value = 14324.667
$ws.Range("M31").Value = -14029.667

$ws.Range("H34").Value = 11717.195
$ws.Range("I34").Value = 14324.667
$ws.Range("K34").Value = 14324.667
$ws.Range("M34").Value = -14122.667

$ws.Range("H105").Value = 6250718
$ws.Range("I105").Value = 6579671.5
$ws.Range("J105").Value = 600
$ws.Range("K105").Value = 6579671.5
$ws.Range("L105").Value = 600
$ws.Range("M105").Value = -6577924.5
$ws.Range("N105").Value = -4094

$ws.Range("H113").Value = 839.2727
$ws.Range("I113").Value = 484.4
$ws.Range("J113").Value = 1135
$ws.Range("K113").Value = 484.4
$ws.Range("L113").Value = 1135
$ws.Range("M113").Value = 1685.6
$ws.Range("N113").Value = -5475

$ws.Range("H122").Value = 1042.7778
$ws.Range("I122").Value = 897.381
$ws.Range("J122").Value = 1246.3334
$ws.Range("K122").Value = 2692.143
$ws.Range("L122").Value = 3739.0002
$ws.Range("M122").Value = -242.143
$ws.Range("N122").Value = -8639.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1911.8889
$ws.Range("I5").Value = 1487
$ws.Range("J5").Value = 3399
$ws.Range("K5").Value = 4461
$ws.Range("L5").Value = 10197
$ws.Range("M5").Value = -4349
$ws.Range("N5").Value = -10421

$ws.Range("H23").Value = 1611.2222
$ws.Range("J23").Value = 1625.125
$ws.Range("L23").Value = 4875.375
$ws.Range("N23").Value = -5345.375

$ws.Range("H131").Value = 147873.69
$ws.Range("J131").Value = 157085.64
$ws.Range("L131").Value = 471256.92
$ws.Range("N131").Value = -481336.92

$ws.Range("H135").Value = 1911.8889
$ws.Range("I135").Value = 1487
$ws.Range("J135").Value = 3399
$ws.Range("K135").Value = 13383
$ws.Range("L135").Value = 30591
$ws.Range("M135").Value = -10848
$ws.Range("N135").Value = -35661

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 33338000
$ws.Range("J52").Value = 33338000
$ws.Range("L52").Value = 33338000
$ws.Range("N52").Value = -33338518

$ws.Range("H70").Value = 12233
$ws.Range("I70").Value = 19601.334
$ws.Range("J70").Value = 4864.6665
$ws.Range("K70").Value = 19601.334
$ws.Range("L70").Value = 4864.6665
$ws.Range("M70").Value = -19331.334
$ws.Range("N70").Value = -5404.6665

$ws.Range("H73").Value = 12233
$ws.Range("I73").Value = 19601.334
$ws.Range("J73").Value = 4864.6665
$ws.Range("K73").Value = 19601.334
$ws.Range("L73").Value = 4864.6665
$ws.Range("M73").Value = -18665.334
$ws.Range("N73").Value = -6736.6665

$ws.Range("H97").Value = 919.6429000000001
$ws.Range("I97").Value = 870.7692
$ws.Range("K97").Value = 870.7692
$ws.Range("M97").Value = -374.7692

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1558.4667
$ws.Range("I100").Value = 1097.1428
$ws.Range("J100").Value = 1962.125
$ws.Range("K100").Value = 1097.1428
$ws.Range("L100").Value = 1962.125
$ws.Range("M100").Value = -556.1428000000001
$ws.Range("N100").Value = -3044.125

$ws.Range("H122").Value = 893633.3
$ws.Range("I122").Value = 1636269.4
$ws.Range("J122").Value = 2470
$ws.Range("K122").Value = 4908808.199999999
$ws.Range("L122").Value = 7410
$ws.Range("M122").Value = -4906358.199999999
$ws.Range("N122").Value = -12310

$ws.Range("H132").Value = 1967.25
$ws.Range("I132").Value = 1600.7273
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 4802.1819
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -2272.1819
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 21277728
$ws.Range("I136").Value = 25001034
$ws.Range("K136").Value = 75003102
$ws.Range("M136").Value = -75000552

$ws.Range("H137").Value = 46614.5
$ws.Range("I137").Value = 45000
$ws.Range("J137").Value = 47306.43
$ws.Range("K137").Value = 45000
$ws.Range("L137").Value = 47306.43
$ws.Range("M137").Value = -39900
$ws.Range("N137").Value = -57506.43
